$d = $word.ActiveDocument

$d.Content.Find.ClearFormatting()
$d.Content.Find.Execute("<id>p097v_1</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p097v_1</id>", 2)

$d.Content.Find.ClearFormatting()
$d.Content.Find.Execute("<id>p097v_2</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p097v_2</id>", 2)
